$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.016.99'
$ws.Range("E2").Value = '  -0.33%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.828.31'
$ws.Range("E3").Value = '  +0.21%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("E4").Value = '  -0.54%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.57'
$ws.Range("E5").Value = '  +0.11%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.006'
$ws.Range("E6").Value = '  -0.35%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4592'
$ws.Range("E7").Value = '  -0.82%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3704'
$ws.Range("E8").Value = '  +2.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07328'
$ws.Range("E9").Value = '  +0.61%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8744'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07937'
$ws.Range("E11").Value = '  +3.97%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '19.79'
$ws.Range("E12").Value = '  -1.59%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.893.51'
$ws.Range("E13").Value = '  +3.07%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.333'
$ws.Range("E14").Value = '  -0.11%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.554'
$ws.Range("E15").Value = '  +1.31%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.61'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008897'
$ws.Range("E18").Value = '  +3.41%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.005'
$ws.Range("E19").Value = '  -0.42%  '

$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.69'
$ws.Range("E20").Value = '  +1.50%  '

$ws.Range("B21").Value = 'WrappedBTC'
$ws.Range("C21").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '26.932.83'
$ws.Range("E21").Value = '  -1.87%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.100'
$ws.Range("E22").Value = '  -2.09%  '

$ws.Range("E23").Value = '  -0.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.152.39'
$ws.Range("E24").Value = '  +2.55%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.01'
$ws.Range("E25").Value = '  +1.30%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.846'
$ws.Range("E26").Value = '  -1.91%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.39'
$ws.Range("E27").Value = '  +0.88%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.051'
$ws.Range("E28").Value = '  -1.61%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.141'
$ws.Range("E29").Value = '  +0.62%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.23'
$ws.Range("E30").Value = '  -0.76%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.958'
$ws.Range("E32").Value = '  +0.00%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7312'
$ws.Range("E33").Value = '  -0.96%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.444'
$ws.Range("E34").Value = '  -0.16%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.133'
$ws.Range("E35").Value = '  -1.29%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.074'
$ws.Range("E36").Value = '  -0.72%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05233'
$ws.Range("E37").Value = '  -0.04%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.430'
$ws.Range("E38").Value = '  -2.77%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01939'
$ws.Range("E39").Value = '  +1.29%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.943'
$ws.Range("E40").Value = '  +0.45%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.128'
$ws.Range("E41").Value = '  -0.40%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5149'
$ws.Range("E42").Value = '  -0.91%  '

$ws.Range("E43").Value = '  +0.23%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.217'
$ws.Range("E44").Value = '  -0.84%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4828'
$ws.Range("E45").Value = '  -0.33%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.22'
$ws.Range("E46").Value = '  +0.72%  '

$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.006'
$ws.Range("E47").Value = '  -0.37%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '102.24'
$ws.Range("E48").Value = '  -1.30%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.626'
$ws.Range("E49").Value = '  -0.54%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06210'
$ws.Range("E50").Value = '  -0.92%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '64.45'
$ws.Range("E51").Value = '  -0.14%  '
